# Fruta / hortaliza, semanal
#
# The underlying data rows (2-20) were re-sorted/shuffled; row 1 (header)
# and row 21 are unaffected. Capture every source data row (columns A:R)
# first, then write them back out in the new row order, so the permutation
# (which includes one long 18-cycle) never reads an already-overwritten row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (values currently sitting in sourceRow, before
# any writes, must end up in destinationRow)
$mapping = @{
    2  = 4
    3  = 19
    4  = 20
    5  = 16
    6  = 7
    7  = 17
    8  = 2
    9  = 3
    10 = 5
    11 = 6
    12 = 11
    13 = 12
    14 = 13
    15 = 15
    16 = 14
    17 = 8
    18 = 9
    19 = 10
    20 = 18
}

# Snapshot every source row's A:R values before writing anything back.
$snapshots = @{}
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    if (-not $snapshots.ContainsKey($srcRow)) {
        $snapshots[$srcRow] = $ws.Range("A$srcRow`:R$srcRow").Value2()
    }
}

# Now write each snapshot into its destination row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $ws.Range("A$destRow`:R$destRow").Value2 = $snapshots[$srcRow]
}
